$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert a new table (ID / City / Country / Continent) before the existing
#    table, followed by an empty paragraph that separates the two tables.
# ---------------------------------------------------------------------------

# Insert an empty paragraph right before the existing table first; the new
# table will then be created in the space before that paragraph, giving the
# final order [new table][blank paragraph][existing table].
$startRng = $d.Range(0, 0)
$startRng.InsertParagraphBefore()

$tableInsertRng = $d.Range(0, 0)
$newTable = $d.Tables.Add($tableInsertRng, 3, 4)
$newTable.Style = "Table Grid"
$newTable.ApplyStyleHeadingRows = $true
$newTable.ApplyStyleLastRow = $false
$newTable.ApplyStyleFirstColumn = $true
$newTable.ApplyStyleLastColumn = $false
$newTable.ApplyStyleRowBands = $true
$newTable.ApplyStyleColumnBands = $false

# Column widths (dxa/twips -> points expected by the Width property).
for ($r = 1; $r -le 3; $r++) {
    $newTable.Cell($r, 1).Width = 2096 / 20.0
    $newTable.Cell($r, 2).Width = 2216 / 20.0
    $newTable.Cell($r, 3).Width = 2142 / 20.0
    $newTable.Cell($r, 4).Width = 2040 / 20.0
}

# Header row text.
$newTable.Cell(1, 1).Range.Text = "ID"
$newTable.Cell(1, 2).Range.Text = "City"
$newTable.Cell(1, 3).Range.Text = "Country"
$newTable.Cell(1, 4).Range.Text = "Continent"
$newTable.Cell(1, 1).Range.ParagraphFormat.Alignment = 1

# Data rows.
$newTable.Cell(2, 1).Range.Text = "0"
$newTable.Cell(2, 2).Range.Text = "Curitiba"
$newTable.Cell(2, 3).Range.Text = "Brazil"
$newTable.Cell(2, 4).Range.Text = "America"

$newTable.Cell(3, 1).Range.Text = "1"
$newTable.Cell(3, 2).Range.Text = "New York"
$newTable.Cell(3, 3).Range.Text = "USA"
$newTable.Cell(3, 4).Range.Text = "America"

# ---------------------------------------------------------------------------
# 2. Fix up the original table (now Tables(2)).
# ---------------------------------------------------------------------------
$origTable = $d.Tables(2)

# 2a. "N" + "ei" + "me" runs -> single run "Neime" (text unchanged, runs
#     merged). A same-value assignment is a no-op in this engine, so we
#     round-trip through a placeholder value first to force the rebuild.
$nameCell = $origTable.Cell(1, 2)
$nameRng = $d.Range($nameCell.Range.Start, $nameCell.Range.End - 1)
$nameRng.Text = "Neime_tmp"
$nameRng2 = $d.Range($nameCell.Range.Start, $nameCell.Range.End - 1)
$nameRng2.Text = "Neime"

# 2b. "Age" run -> "Ag" + "i" runs (two runs, same visible text "Agi").
$ageCell = $origTable.Cell(1, 3)
$ageRng = $d.Range($ageCell.Range.Start, $ageCell.Range.End - 1)
$ageRng.Text = "Agi_tmp"
$ageRng2 = $d.Range($ageCell.Range.Start, $ageCell.Range.End - 1)
$ageRng2.Text = "Agi"
# Force a run split right before the last character by toggling a
# character property on/off (identical before/after formatting).
$ageCellEnd = $ageCell.Range.End - 1
$splitPoint = $d.Range($ageCellEnd - 1, $ageCellEnd)
$splitPoint.Bold = 1
$splitPoint.Bold = 0

# ---------------------------------------------------------------------------
# 3. Add a new row (2, Bob, 56) to the original table.
# ---------------------------------------------------------------------------
$newRow = $origTable.Rows.Add()
$newRowIndex = $origTable.Rows.Count
$origTable.Cell($newRowIndex, 1).Range.Text = "2"
$origTable.Cell($newRowIndex, 2).Range.Text = "Bob"
$origTable.Cell($newRowIndex, 3).Range.Text = "56"

Write-Output "edit complete"
